$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.262.34'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.11%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.543.48'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +5.15%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '527.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.27%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.74'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.56%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("E8").Value = '  +3.38%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.542.24'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.84%  '

$ws.Range("E10").Value = '  +4.40%  '

$ws.Range("E11").Value = '  -0.86%  '

$ws.Range("E12").Value = '  +0.57%  '

$ws.Range("E13").Value = '  +1.96%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.997.17'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '59.177.59'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.10%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.44'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.03%  '

$ws.Range("E17").Value = '  +3.84%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.546.91'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.88%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.74'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.98%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '323.97'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.17%  '

$ws.Range("E21").Value = '  +3.35%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.17'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +10.00%  '

$ws.Range("E23").Value = '  +0.27%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.46'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.97%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.411'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.78%  '

$ws.Range("E26").Value = '  +0.07%  '

$ws.Range("E27").Value = '  +1.49%  '

$ws.Range("E28").Value = '  +4.71%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0760'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.48%  '

$ws.Range("E30").Value = '  +7.69%  '

$ws.Range("E31").Value = '  +4.46%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '169.34'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.05%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.37'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.56%  '

$ws.Range("E34").Value = '  -0.01%  '

$ws.Range("E35").Value = '  -0.12%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.29'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.76%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.27'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.12%  '

$ws.Range("E38").Value = '  +3.84%  '

$ws.Range("E39").Value = '  +5.95%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.80'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.27%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.787'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.53%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '280.54'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.19%  '

$ws.Range("E43").Value = '  +10.90%  '

$ws.Range("E44").Value = '  +4.37%  '

$ws.Range("E45").Value = '  +5.09%  '

$ws.Range("E46").Value = '  +3.72%  '

$ws.Range("E47").Value = '  +2.87%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0508'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.97%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '17.93'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.52%  '

$ws.Range("E50").Value = '  +4.39%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.17'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.39%  '
